$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C.
# This shifts the former C..G (Régime, Type, Nombre de personne, Instructions)
# to D..H.
$ws.Columns.Item(3).Insert()

# New column header: "Quantité Ingredient" (numeric quantity column,
# complementing the existing "Quantité Ingrédients" text column that is
# now in column D).
$ws.Cells.Item(1, 3).Value = "Quantité Ingredient"

# Give the new column the same width as column B.
$ws.Range("C1").ColumnWidth = 21.83

# Fill in the numeric quantities for each ingredient row.
$ws.Cells.Item(2, 3).Value = 4
$ws.Cells.Item(3, 3).Value = 3
$ws.Cells.Item(4, 3).Value = 4
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(6, 3).Value = 0.5
$ws.Cells.Item(7, 3).Value = 3

# Match the new active selection.
$ws.Range("D10").Select()
